# Auto-generated edit script: apply 2024-10-18 data update
# Updates year-2024 (column K) running totals plus a couple of revised
# 2015 (column B) historical values, across the Citywide Totals sheet,
# the By Neighborhood summary sheet, and each affected neighborhood sheet.

$wb = $excel.ActiveWorkbook

$updates = [ordered]@{
    'Citywide Totals' = [ordered]@{ 'K2' = 6471; 'K3' = 6666; 'B4' = 1705; 'K4' = 1390; 'K5' = 478; 'K6' = 7341; 'B7' = 23338; 'K7' = 22346 }
    'Logan Square' = [ordered]@{ 'K6' = 120; 'K7' = 284 }
    'Austin' = [ordered]@{ 'K3' = 446; 'K7' = 1466 }
    'South Chicago' = [ordered]@{ 'K6' = 111; 'K7' = 480 }
    'West Pullman' = [ordered]@{ 'K2' = 126; 'K7' = 367 }
    'Grand Crossing' = [ordered]@{ 'K2' = 217; 'K4' = 36; 'K7' = 759 }
    'New City' = [ordered]@{ 'K6' = 189; 'K7' = 523 }
    'Woodlawn' = [ordered]@{ 'K3' = 153; 'K7' = 371 }
    'By Neighborhood' = [ordered]@{ 'K2' = 196; 'K4' = 82; 'K7' = 671; 'K8' = 1466; 'K11' = 412; 'K14' = 113; 'K15' = 230; 'K16' = 57; 'K19' = 654; 'B22' = 63; 'K26' = 31; 'K29' = 1207; 'K34' = 127; 'K36' = 284; 'K37' = 759; 'K42' = 827; 'K43' = 183; 'K48' = 281; 'K50' = 106; 'K52' = 587; 'K53' = 284; 'K60' = 131; 'K63' = 63; 'K65' = 523; 'K67' = 876; 'K72' = 117; 'K73' = 199; 'K83' = 480; 'K84' = 180; 'K86' = 136; 'K89' = 331; 'K90' = 209; 'K94' = 299; 'K95' = 367; 'K96' = 239; 'K99' = 371; 'B101' = 23338; 'K101' = 22346 }
    'North Lawndale' = [ordered]@{ 'K2' = 241; 'K3' = 319; 'K6' = 247; 'K7' = 876 }
    'South Deering' = [ordered]@{ 'K3' = 72; 'K7' = 180 }
    'Englewood' = [ordered]@{ 'K2' = 344; 'K6' = 348; 'K7' = 1207 }
    'Lake View' = [ordered]@{ 'K2' = 43; 'K6' = 131; 'K7' = 281 }
    'Chatham' = [ordered]@{ 'K2' = 193; 'K3' = 197; 'K7' = 654 }
    'Bridgeport' = [ordered]@{ 'K3' = 26; 'K6' = 41; 'K7' = 113 }
    'Humboldt Park' = [ordered]@{ 'K2' = 223; 'K3' = 251; 'K5' = 13; 'K6' = 307; 'K7' = 827 }
    'West Ridge' = [ordered]@{ 'K2' = 75; 'K7' = 239 }
    'Grand Boulevard' = [ordered]@{ 'K2' = 111; 'K7' = 284 }
    'Auburn Gresham' = [ordered]@{ 'K2' = 220; 'K3' = 221; 'K5' = 27; 'K7' = 671 }
    'Garfield Ridge' = [ordered]@{ 'K2' = 50; 'K7' = 127 }
    'West Loop' = [ordered]@{ 'K3' = 62; 'K6' = 135; 'K7' = 299 }
    'Brighton Park' = [ordered]@{ 'K6' = 71; 'K7' = 230 }
    'Lincoln Square' = [ordered]@{ 'K2' = 28; 'K7' = 106 }
    'East Village' = [ordered]@{ 'K3' = 6; 'K7' = 31 }
    'Belmont Cragin' = [ordered]@{ 'K4' = 23; 'K7' = 412 }
    'Portage Park' = [ordered]@{ 'K2' = 66; 'K6' = 68; 'K7' = 199 }
    'Albany Park' = [ordered]@{ 'K6' = 62; 'K7' = 196 }
    'Uptown' = [ordered]@{ 'K2' = 95; 'K7' = 331 }
    'Streeterville' = [ordered]@{ 'K3' = 23; 'K6' = 33; 'K7' = 136 }
    'Washington Heights' = [ordered]@{ 'K3' = 61; 'K7' = 209 }
    'Morgan Park' = [ordered]@{ 'K3' = 39; 'K7' = 131 }
    'Hyde Park' = [ordered]@{ 'K3' = 49; 'K4' = 25; 'K7' = 183 }
    'Clearing' = [ordered]@{ 'B4' = 9; 'B7' = 63 }
    'Old Town' = [ordered]@{ 'K2' = 26; 'K3' = 28; 'K7' = 117 }
    'Little Village' = [ordered]@{ 'K2' = 154; 'K6' = 213; 'K7' = 587 }
    'Archer Heights' = [ordered]@{ 'K2' = 27; 'K7' = 82 }
    'Bucktown' = [ordered]@{ 'K3' = 6; 'K7' = 57 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
